$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.011751144388498
$ws.Range("D2").Value = 1.014193705388803
$ws.Range("E2").Value = 1.013890478250391
$ws.Range("F2").Value = 1.014488441360307
$ws.Range("I2").Value = 1.023013711196102
$ws.Range("J2").Value = 1.016997814213541
$ws.Range("K2").Value = 1.01705259004418
$ws.Range("L2").Value = 1.016750272417264
$ws.Range("M2").Value = 1.017346442540627
$ws.Range("N2").Value = 1.018442067761004

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.013292748937411
$ws.Range("D3").Value = 1.015558820689204
$ws.Range("E3").Value = 1.015209691036573
$ws.Range("F3").Value = 1.017253898378762
$ws.Range("I3").Value = 1.023304341867529
$ws.Range("J3").Value = 1.018169665542221
$ws.Range("K3").Value = 1.018221041842741
$ws.Range("L3").Value = 1.017872882253083
$ws.Range("M3").Value = 1.019911419986296
$ws.Range("N3").Value = 1.019615583252984

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.014279679616375
$ws.Range("D4").Value = 1.016432694081755
$ws.Range("E4").Value = 1.01605456773879
$ws.Range("F4").Value = 1.018992688228816
$ws.Range("I4").Value = 1.023479063045243
$ws.Range("J4").Value = 1.018917109664778
$ws.Range("K4").Value = 1.018967170918123
$ws.Range("L4").Value = 1.018590041386176
$ws.Range("M4").Value = 1.021520437463888
$ws.Range("N4").Value = 1.020364088831928

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.014692095914894
$ws.Range("D5").Value = 1.01679784959094
$ws.Range("E5").Value = 1.016407699411662
$ws.Range("F5").Value = 1.019711738599856
$ws.Range("I5").Value = 1.023549349751155
$ws.Range("J5").Value = 1.019228785847496
$ws.Range("K5").Value = 1.019278503880662
$ws.Range("L5").Value = 1.018889358919755
$ws.Range("M5").Value = 1.022184911517033
$ws.Range("N5").Value = 1.020676207630576

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.014761197737103
$ws.Range("D6").Value = 1.016859031717523
$ws.Range("E6").Value = 1.016466872370215
$ws.Range("F6").Value = 1.019831775563913
$ws.Range("I6").Value = 1.023560966267709
$ws.Range("J6").Value = 1.019280969313923
$ws.Range("K6").Value = 1.019330641951534
$ws.Range("L6").Value = 1.018939489071833
$ws.Range("M6").Value = 1.022295783464792
$ws.Range("N6").Value = 1.02072846520351

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.014285200064878
$ws.Range("D7").Value = 1.016437581981028
$ws.Range("E7").Value = 1.016059294324445
$ws.Range("F7").Value = 1.019002342877297
$ws.Range("I7").Value = 1.023480014626097
$ws.Range("J7").Value = 1.01892128426061
$ws.Range("K7").Value = 1.01897134010692
$ws.Range("L7").Value = 1.018594049384661
$ws.Range("M7").Value = 1.021529362950481
$ws.Range("N7").Value = 1.020368269356165

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.012274358294898
$ws.Range("D8").Value = 1.014657034071791
$ws.Range("E8").Value = 1.014338146382617
$ws.Range("F8").Value = 1.015433644118661
$ws.Range("I8").Value = 1.023114709385951
$ws.Range("J8").Value = 1.017396113452992
$ws.Range("K8").Value = 1.017449555860714
$ws.Range("L8").Value = 1.017131600805361
$ws.Range("M8").Value = 1.018223897262547
$ws.Range("N8").Value = 1.018840932631061

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.00864767681999
$ws.Range("D9").Value = 1.01144515747268
$ws.Range("E9").Value = 1.011236461546351
$ws.Range("F9").Value = 1.008748158077022
$ws.Range("I9").Value = 1.022367526740987
$ws.Range("J9").Value = 1.01462368178593
$ws.Range("K9").Value = 1.014689973466053
$ws.Range("L9").Value = 1.0144819953177
$ws.Range("M9").Value = 1.012002275170437
$ws.Range("N9").Value = 1.016064563793028

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.006170559489222
$ws.Range("D10").Value = 1.009250984001498
$ws.Range("E10").Value = 1.009119641440275
$ws.Range("F10").Value = 1.004010528187619
$ws.Range("I10").Value = 1.021797782711332
$ws.Range("J10").Value = 1.01271531083376
$ws.Range("K10").Value = 1.012794957402486
$ws.Range("L10").Value = 1.01266411213005
$ws.Range("M10").Value = 1.007574454944737
$ws.Range("N10").Value = 1.014153482735213

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.005083065602942
$ws.Range("D11").Value = 1.008287612640525
$ws.Range("E11").Value = 1.008190731961095
$ws.Range("F11").Value = 1.001889156436911
$ws.Range("I11").Value = 1.021533572245972
$ws.Range("J11").Value = 1.011873980239429
$ws.Range("K11").Value = 1.011960585483849
$ws.Range("L11").Value = 1.011864086257473
$ws.Range("M11").Value = 1.005587498512735
$ws.Range("N11").Value = 1.013310957354936

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.004676811853962
$ws.Range("D12").Value = 1.00792771249389
$ws.Range("E12").Value = 1.007843781866891
$ws.Range("F12").Value = 1.001090365097118
$ws.Range("I12").Value = 1.021432753711119
$ws.Range("J12").Value = 1.011559152372105
$ws.Range("K12").Value = 1.011648522776114
$ws.Range("L12").Value = 1.011564927338318
$ws.Range("M12").Value = 1.004838686493379
$ws.Range("N12").Value = 1.01299568239593

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.004764060433098
$ws.Range("D13").Value = 1.008005006673628
$ws.Range("E13").Value = 1.007918291389118
$ws.Range("F13").Value = 1.001262203374373
$ws.Range("I13").Value = 1.021454501683428
$ws.Range("J13").Value = 1.011626790148051
$ws.Range("K13").Value = 1.011715559189739
$ws.Range("L13").Value = 1.0116291891676
$ws.Range("M13").Value = 1.004999801515209
$ws.Range("N13").Value = 1.013063416225276

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.005049532108908
$ws.Range("D14").Value = 1.008257905624941
$ws.Range("E14").Value = 1.008162092329393
$ws.Range("F14").Value = 1.001823350679566
$ws.Range("I14").Value = 1.021525293500321
$ws.Range("J14").Value = 1.011848004208179
$ws.Range("K14").Value = 1.011934834318987
$ws.Range("L14").Value = 1.011839398719397
$ws.Range("M14").Value = 1.005525822942547
$ws.Range("N14").Value = 1.013284944434743

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.005225112237996
$ws.Range("D15").Value = 1.008413449851108
$ws.Range("E15").Value = 1.008312050889605
$ws.Range("F15").Value = 1.002167648726334
$ws.Range("I15").Value = 1.021568554160169
$ws.Range("J15").Value = 1.011983991942022
$ws.Range("K15").Value = 1.01206965143403
$ws.Range("L15").Value = 1.011968649717498
$ws.Range("M15").Value = 1.005848485954044
$ws.Range("N15").Value = 1.013421125286766

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.006242413066293
$ws.Range("D16").Value = 1.009314634458733
$ws.Range("E16").Value = 1.009181025496012
$ws.Range("F16").Value = 1.004149816903949
$ws.Range("I16").Value = 1.021814944556084
$ws.Range("J16").Value = 1.012770825607459
$ws.Range("K16").Value = 1.012850035439281
$ws.Range("L16").Value = 1.012716931008891
$ws.Range("M16").Value = 1.007704829532407
$ws.Range("N16").Value = 1.014209076346259

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.006876503078135
$ws.Range("D17").Value = 1.00987632310218
$ws.Range("E17").Value = 1.009722771417489
$ws.Range("F17").Value = 1.00537423912678
$ws.Range("I17").Value = 1.021964778288326
$ws.Range("J17").Value = 1.013260325043385
$ws.Range("K17").Value = 1.013335806549797
$ws.Range("L17").Value = 1.013182821120296
$ws.Range("M17").Value = 1.008850404445092
$ws.Range("N17").Value = 1.01469927092753

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.007244926018013
$ws.Range("D18").Value = 1.010202669960987
$ws.Range("E18").Value = 1.010037579073759
$ws.Range("F18").Value = 1.006081697157302
$ws.Range("I18").Value = 1.022050488041301
$ws.Range("J18").Value = 1.013544399941989
$ws.Range("K18").Value = 1.013617819972301
$ws.Range("L18").Value = 1.013453329389758
$ws.Range("M18").Value = 1.009511897001076
$ws.Range("N18").Value = 1.014983749245075

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.00737030815299
$ws.Range("D19").Value = 1.010313731203653
$ws.Range("E19").Value = 1.010144721411059
$ws.Range("F19").Value = 1.006321790182017
$ws.Range("I19").Value = 1.022079428339569
$ws.Range("J19").Value = 1.013641019561814
$ws.Range("K19").Value = 1.013713755772717
$ws.Range("L19").Value = 1.013545357568964
$ws.Range("M19").Value = 1.009736320855286
$ws.Range("N19").Value = 1.015080506075842

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.006808619757889
$ws.Range("D20").Value = 1.009816191737622
$ws.Range("E20").Value = 1.009664770065313
$ws.Range("F20").Value = 1.005243568167238
$ws.Range("I20").Value = 1.021948877249654
$ws.Range("J20").Value = 1.013207955998
$ws.Range("K20").Value = 1.013283825784255
$ws.Range("L20").Value = 1.013132963948074
$ws.Range("M20").Value = 1.008728190405761
$ws.Range("N20").Value = 1.014646827512094

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.004965532200604
$ws.Range("D21").Value = 1.008183490635576
$ws.Range("E21").Value = 1.008090352320266
$ws.Range("F21").Value = 1.001658408055019
$ws.Range("I21").Value = 1.021504521438549
$ws.Range("J21").Value = 1.011782926796021
$ws.Range("K21").Value = 1.01187032286249
$ws.Range("L21").Value = 1.011777552783077
$ws.Range("M21").Value = 1.005371222496674
$ws.Range("N21").Value = 1.013219774605195

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.003793304507564
$ws.Range("D22").Value = 1.007144985821511
$ws.Range("E22").Value = 1.007089358998104
$ws.Range("F22").Value = 0.9993415083748991
$ws.Range("I22").Value = 1.021209610501494
$ws.Range("J22").Value = 1.010873493545802
$ws.Range("K22").Value = 1.010969182411109
$ws.Range("L22").Value = 1.01091378453477
$ws.Range("M22").Value = 1.003198100317083
$ws.Range("N22").Value = 1.012309049855447

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.0044160210813
$ws.Range("D23").Value = 1.00769667376586
$ws.Range("E23").Value = 1.007621077844902
$ws.Range("F23").Value = 1.000575799664811
$ws.Range("I23").Value = 1.021367437394189
$ws.Range("J23").Value = 1.011356900995297
$ws.Range("K23").Value = 1.011448093248226
$ws.Range("L23").Value = 1.011372801958976
$ws.Range("M23").Value = 1.0043561396138
$ws.Range("N23").Value = 1.012793143798968

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.006839297725893
$ws.Range("D24").Value = 1.009843366453837
$ws.Range("E24").Value = 1.009690982033348
$ws.Range("F24").Value = 1.005302633514715
$ws.Range("I24").Value = 1.021956067450449
$ws.Range("J24").Value = 1.013231623766641
$ws.Range("K24").Value = 1.01330731775643
$ws.Range("L24").Value = 1.013155496078587
$ws.Range("M24").Value = 1.008783434359067
$ws.Range("N24").Value = 1.014670528891682

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.009595462129415
$ws.Range("D25").Value = 1.012284603828029
$ws.Range("E25").Value = 1.01204674862512
$ws.Range("F25").Value = 1.010524912778679
$ws.Range("I25").Value = 1.022573152998984
$ws.Range("J25").Value = 1.015350778269024
$ws.Range("K25").Value = 1.015412916725442
$ws.Range("L25").Value = 1.015175847790463
$ws.Range("M25").Value = 1.013659052224085
$ws.Range("N25").Value = 1.016792692836528
